$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 365, shifting existing rows 365:379 down to 366:380
$ws.Rows("365:365").Insert()

# Populate the newly inserted row 365 with its data
$ws.Range("A365").Value = 7
$ws.Range("B365").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C365").Value = "Ñuble"
$ws.Range("D365").Value = 44516
$ws.Range("E365").Value = 16
$ws.Range("F365").Value = 100112004
$ws.Range("G365").Value = "Cebolla"
$ws.Range("H365").Value = "Sin especificar"
$ws.Range("I365").Value = "1a nueva(o)"
$ws.Range("J365").Value = 16000
$ws.Range("K365").Value = 800
$ws.Range("L365").Value = 900
$ws.Range("M365").Value = 850
$ws.Range("N365").Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O365").Value = "Región del Maule"
$ws.Range("P365").Value = 85
$ws.Range("Q365").Value = 10
$ws.Range("R365").Value = "Hortaliza"
